# This workbook originally has two sheets, in this order:
#   1) hotel_info   (sheetId=1)  -> 9 columns A..I, one header row + one data row
#   2) review_info  (sheetId=2)  -> 25 columns A..Y, header row only
#
# The target state:
#   - sheet order on the tab strip becomes: review_info (first), hotel_info (second)
#   - hotel_info gains a new "State" column right after "Hotel_Name" (value "Louisiana"
#     for the only data row, the Luling / Louisiana Motel 6)
#
# The automation surface here cannot reorder worksheets (Worksheet.Move /
# Worksheets.Add(Before:=/After:=) are accepted but not implemented by the
# engine), and deleting + re-adding sheets mints brand new sheetId values
# (they are never reused), which would not match the target's sheetId
# numbering. So instead of moving sheets, we keep the two original sheet
# objects exactly where they are and simply swap what is written into them:
#   - the sheet currently in tab position 1 (hotel_info, sheetId=1) becomes
#     "review_info" and gets the review_info header row
#   - the sheet currently in tab position 2 (review_info, sheetId=2) becomes
#     "hotel_info" and gets the (new, wider) hotel_info header + data row
# That reproduces the desired final tab order/names/ids without needing an
# actual "move".

$wb = $excel.ActiveWorkbook

$sheetA = $wb.Worksheets.Item(1)   # currently "hotel_info"
$sheetB = $wb.Worksheets.Item(2)   # currently "review_info"

# Rename both out of the way first so the final names (which are simply
# swapped between the two sheets) never collide mid-script.
$sheetA.Name = "__tmp_sheetA__"
$sheetB.Name = "__tmp_sheetB__"

# --- Rebuild sheet A as the new "review_info" (header row only) -----------

$sheetA.Cells.Clear()

$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $sheetA.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

$sheetA.Name = "review_info"

# --- Rebuild sheet B as the new "hotel_info" (header row + one data row) --

$sheetB.Cells.Clear()

$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)

for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $sheetB.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

# Row 2 data (A..J). Numeric cells stay numeric; "21"/"1"/"23" stay text
# (leading apostrophe forces text so they keep matching the original
# string-typed cells instead of becoming numbers).
$sheetB.Cells.Item(2, 1).Value = 39114
$sheetB.Cells.Item(2, 2).Value = "Motel 6 Luling"
$sheetB.Cells.Item(2, 3).Value = "Louisiana"
$sheetB.Cells.Item(2, 4).Value = "Luling"
$sheetB.Cells.Item(2, 5).Value = 70070
$sheetB.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g40290-d116272-Reviews-Motel_6_Luling_LA-Luling_Louisiana.html"
$sheetB.Cells.Item(2, 7).Value = "Motel 6 Luling, LA"
$sheetB.Cells.Item(2, 8).Value = "'21"
$sheetB.Cells.Item(2, 9).Value = "'1"
$sheetB.Cells.Item(2, 10).Value = "'23"

$sheetB.Name = "hotel_info"
